$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Choice"
$ws.Range("B1").Value = "Input"
$ws.Range("C1").Value = "Output"
$ws.Range("D1").Value = "Exit loop?"

# Column A - Choice values (establishes shared-string order: Deposit, Withdrawel, View, Exit)
$ws.Range("A2").Value = "Deposit"
$ws.Range("A3").Value = "Deposit"
$ws.Range("A4").Value = "Withdrawel"
$ws.Range("A5").Value = "Withdrawel"
$ws.Range("A6").Value = "View"
$ws.Range("A7").Value = "Exit"

# Column D - Exit loop? values (establishes shared-string order: no, yes)
$ws.Range("D2").Value = "no"
$ws.Range("D3").Value = "no"
$ws.Range("D4").Value = "no"
$ws.Range("D5").Value = "no"
$ws.Range("D6").Value = "no"
$ws.Range("D7").Value = "yes"

# Column B - Input values
$ws.Range("B2").Value = 300
$ws.Range("B3").Value = 45
$ws.Range("B4").Value = 1000
$ws.Range("B5").Value = 315

# Column C - Output formulas / values
$ws.Range("C2").Formula = "=1000 + B2"
$ws.Range("C3").Formula = "=C2+B3"
$ws.Range("C4").Formula = "=C3-B4"
$ws.Range("C5").Formula = "=C4-B5"
$ws.Range("C6").Value = 30

# Selection matches the final state in the diff
$ws.Range("D7").Select()
